$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was date 44489/Primera -> now 44167/Segunda)
$ws.Range("D2").Value = 44167
$ws.Range("L2").Value = 'Segunda'
$ws.Range("N2").Value = 18000
$ws.Range("O2").Value = 19000
$ws.Range("P2").Value = 18500
$ws.Range("Q2").Value = '$/caja 13 kilos'
$ws.Range("S2").Value = 1423
$ws.Range("T2").Value = 13

# Row 3 (was date 44167/Segunda -> now 44489/Primera)
$ws.Range("D3").Value = 44489
$ws.Range("L3").Value = 'Primera'
$ws.Range("N3").Value = 24000
$ws.Range("O3").Value = 25000
$ws.Range("P3").Value = 24500
$ws.Range("Q3").Value = '$/caja 12 kilos'
$ws.Range("S3").Value = 2042
$ws.Range("T3").Value = 12

# Row 4 (was date 44441/Primera -> now 44475/Especial)
$ws.Range("D4").Value = 44475
$ws.Range("L4").Value = 'Especial'
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 32000
$ws.Range("O4").Value = 33000
$ws.Range("P4").Value = 32500
$ws.Range("S4").Value = 2708

# Row 5 (was date 44545 -> now 44468)
$ws.Range("D5").Value = 44468
$ws.Range("N5").Value = 29000
$ws.Range("O5").Value = 30000
$ws.Range("P5").Value = 29500
$ws.Range("Q5").Value = '$/bandeja 10 kilos'
$ws.Range("S5").Value = 2950
$ws.Range("T5").Value = 10

# Row 6 (was date 44475/Especial -> now 44441/Primera)
$ws.Range("D6").Value = 44441
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 29000
$ws.Range("O6").Value = 30000
$ws.Range("P6").Value = 29500
$ws.Range("S6").Value = 2458

# Row 8 (was date 44524 -> now 44496)
$ws.Range("D8").Value = 44496

# Row 9 (was date 44468 -> now 44482)
$ws.Range("D9").Value = 44482
$ws.Range("M9").Value = 160
$ws.Range("N9").Value = 25000
$ws.Range("O9").Value = 26000
$ws.Range("P9").Value = 25500
$ws.Range("Q9").Value = '$/caja 12 kilos'
$ws.Range("S9").Value = 2125
$ws.Range("T9").Value = 12

# Row 10 (was date 44482 -> now 44545)
$ws.Range("D10").Value = 44545
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 23000
$ws.Range("O10").Value = 24000
$ws.Range("P10").Value = 23500
$ws.Range("Q10").Value = '$/bandeja 12 kilos'
$ws.Range("S10").Value = 1958

# Row 11 (was date 44496 -> now 44524)
$ws.Range("D11").Value = 44524
